$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the statistics table (columns B:D, rows 2-10) ---
# Values are written as text (quote-prefixed) so they keep the same shared-string
# storage the workbook already used, then formats are reset so no new per-cell
# style is introduced.

$values = @{
    "B2" = "3986.1732101616626";  "C2" = "472.0692609699769";   "D2" = "2776.6209699769056"
    "B3" = "50.0";                 "C3" = "2.97";                 "D3" = "17.42"
    "B4" = "28327.395400161044";  "C4" = "3228.4989260273464";  "D4" = "19016.09383693793"
    "B5" = "1";                    "C5" = "0.01";                 "D5" = "0.05"
    "B6" = "472183";               "C6" = "53026.45";             "D6" = "312878.52"
    "B7" = "1726013";              "C7" = "204405.99";            "D7" = "1202276.8800000001"
    "B8" = "1";                    "C8" = "0";                    "D8" = "0"
    "B9" = "4.0";                  "C9" = "0.36";                 "D9" = "2.1"
    "B10" = "536.0";               "C10" = "34.01";               "D10" = "199.5"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = "'" + $values[$addr]
}

[void]$ws.Range("B2:D10").ClearFormats()

# --- Column widths (closest values reachable through the pixel-snapped
# ColumnWidth property; targets are 19.28515625 / 17 / 31.85546875 / 32.7109375) ---
$ws.Columns("A").ColumnWidth = 18.5
$ws.Columns("B").ColumnWidth = 16.16666667
$ws.Columns("C").ColumnWidth = 31
$ws.Columns("D").ColumnWidth = 31.83333333

# --- Selection ---
[void]$ws.Range("D28").Select()
